$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.786.89"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").Value = "2.211.46"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'261.17"
$ws.Range("E5").Value = "  +2.42%  "

$ws.Range("D6").Value = "'86.61"
$ws.Range("E6").Value = "  +13.66%  "

$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +1.15%  "

$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("D10").Value = "'45.11"
$ws.Range("E10").Value = "  +8.03%  "

$ws.Range("D11").Value = "'0.0914"
$ws.Range("E11").Value = "  +0.60%  "

$ws.Range("D12").Value = "'7.42"
$ws.Range("E12").Value = "  +7.20%  "

$ws.Range("E13").Value = "  +1.37%  "

$ws.Range("D14").Value = "2.544.81"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").Value = "2.213.28"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").Value = "43.732.68"
$ws.Range("E18").Value = "  +2.02%  "

$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").Value = "'5.94"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("E21").Value = "  -2.07%  "

$ws.Range("E22").Value = "  +6.53%  "

$ws.Range("D23").Value = "'230.32"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "'8.87"
$ws.Range("E24").Value = "  -3.51%  "

$ws.Range("E26").Value = "  +5.78%  "

$ws.Range("D27").Value = "'10.65"
$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("D28").Value = "'39.85"
$ws.Range("E28").Value = "  -3.24%  "

$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +3.05%  "

$ws.Range("E30").Value = "  +1.87%  "

$ws.Range("D31").Value = "'174.69"
$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("D33").Value = "'0.0875"
$ws.Range("E33").Value = "  +3.44%  "

$ws.Range("E34").Value = "  +3.86%  "

$ws.Range("D35").Value = "'0.122"
$ws.Range("E35").Value = "  +1.02%  "

$ws.Range("E36").Value = "  +4.68%  "

$ws.Range("D37").Value = "'4.49"
$ws.Range("E37").Value = "  +4.57%  "

$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").Value = "'2.95"
$ws.Range("E39").Value = "  +6.89%  "

$ws.Range("D40").Value = "'12.58"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").Value = "'63.59"
$ws.Range("E42").Value = "  +6.13%  "

$ws.Range("E43").Value = "  +4.31%  "

$ws.Range("D44").Value = "'0.200"
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("D45").Value = "'100.59"
$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("D46").Value = "'0.0978"
$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("D47").Value = "'8.31"
$ws.Range("E47").Value = "  +0.23%  "

$ws.Range("E48").Value = "  +4.48%  "

$ws.Range("E49").Value = "  +1.45%  "

$ws.Range("D50").Value = "'0.447"
$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("E51").Value = "  +4.84%  "
